$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trend_instructions")

# New headers
$ws.Range("K1").Value = "use_as_trend"
$ws.Range("L1").Value = "match_year"

# New row 2 data
$ws.Range("K2").Value = "T"
$ws.Range("L2").Value = 1953

# Column width for column J (bestFit-like width)
$ws.Columns.Item(10).EntireColumn.AutoFit()

# View adjustments
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("J13").Select()
